$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.59"
$ws.Range("E2").Value = "'1.74%"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'27.30"
$ws.Range("E3").Value = "'2.20%"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'4.674"
$ws.Range("E4").Value = "'0.36%"
$ws.Range("G4").Value = "'12"
$ws.Range("E5").Value = "'2.60%"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.666"
$ws.Range("E6").Value = "'0.69%"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'0.8494"
$ws.Range("E7").Value = "'-0.01%"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'0.9215"
$ws.Range("E8").Value = "'1.29%"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.1399"
$ws.Range("E9").Value = "'1.40%"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.04854"
$ws.Range("E10").Value = "'17.45%"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.07087"
$ws.Range("E11").Value = "'1.18%"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.03078"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.09053"
$ws.Range("E13").Value = "'-0.63%"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.001543"
$ws.Range("E14").Value = "'0.87%"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.0006068"
$ws.Range("E15").Value = "'-94.07%"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.006137"
$ws.Range("E16").Value = "'1.36%"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'3.450"
$ws.Range("E17").Value = "'-0.55%"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'3.145"
$ws.Range("E18").Value = "'-0.27%"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'2.163"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("G19").Value = "'12"
$ws.Range("E20").Value = "'3.40%"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.1305"
$ws.Range("E21").Value = "'0.90%"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'4.092"
$ws.Range("E22").Value = "'5.72%"
$ws.Range("G22").Value = "'12"
$ws.Range("D23").Value = "'0.04251"
$ws.Range("E23").Value = "'0.34%"
$ws.Range("G23").Value = "'12"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.82%"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.003797"
$ws.Range("E25").Value = "'-19.91%"
$ws.Range("G25").Value = "'12"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("G26").Value = "'12"
$ws.Range("E27").Value = "'3.40%"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03850"
$ws.Range("E40").Value = "'1.93%"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.1111"
$ws.Range("E41").Value = "'1.71%"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.004085"
$ws.Range("E42").Value = "'-34.23%"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.01633"
$ws.Range("E43").Value = "'17.75%"
$ws.Range("G43").Value = "'12"
$ws.Range("E44").Value = "'0.77%"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005176"
$ws.Range("E45").Value = "'0.29%"
$ws.Range("G45").Value = "'12"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("G46").Value = "'12"
$ws.Range("E47").Value = "'-43.74%"
$ws.Range("G47").Value = "'12"
$ws.Range("E48").Value = "'36.12%"
$ws.Range("G48").Value = "'12"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("G49").Value = "'12"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"